$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the path separator in the default value for fname_sim (row 3, column D)
$ws.Range("D3").Value = "saves/Ns.txt"

# 2. Add a new "debug" setting row at the bottom of the settings table
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "Отладка"
$ws.Range("C45").Value = "debug"
$ws.Range("D45").Value = 1
$ws.Range("F45").Value = "bool"
$ws.Range("G45").Value = "numerical"

# 3. Update the view so the new row is visible / selected
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G45").Select()
